# Update Gdf11-Acvr2b.xlsx values with new TPM-derived statistics (per commit "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("E2").Value = 3
    $ws.Range("F2").Value = 1
    $ws.Range("G2").Value = 1.419475
    $ws.Range("H2").Value = 4.258425
    $ws.Range("I2").Value = 0.1541931834006784
    $ws.Range("J2").Value = 0.1541931834006784
    $ws.Range("M2").Value = 1.716657
    $ws.Range("N2").Value = 5.149971
    $ws.Range("O2").Value = 0.3840886036988016
    $ws.Range("P2").Value = 0.3840886036988015
    $ws.Range("Q2").Value = 2.436751695075
    $ws.Range("R2").Value = 21.930765255675
    $ws.Range("S2").Value = 0.0592238445122398
    $ws.Range("T2").Value = 0.05922384451223979
    # Row 3
    $ws.Range("E3").Value = 3
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 1.419475
    $ws.Range("H3").Value = 4.258425
    $ws.Range("I3").Value = 0.1541931834006784
    $ws.Range("J3").Value = 0.1541931834006784
    $ws.Range("O3").Value = 0.07870146593648156
    $ws.Range("P3").Value = 0.07870146593648154
    $ws.Range("Q3").Value = 0.4993012775666668
    $ws.Range("R3").Value = 4.493711498100001
    $ws.Range("S3").Value = 0.01213522957104615
    $ws.Range("T3").Value = 0.01213522957104614
    # Row 4
    $ws.Range("E4").Value = 3
    $ws.Range("F4").Value = 1
    $ws.Range("G4").Value = 1.419475
    $ws.Range("H4").Value = 4.258425
    $ws.Range("I4").Value = 0.1541931834006784
    $ws.Range("J4").Value = 0.1541931834006784
    $ws.Range("M4").Value = 1.677572333333333
    $ws.Range("N4").Value = 5.032717
    $ws.Range("O4").Value = 0.3753437146230962
    $ws.Range("P4").Value = 0.3753437146230962
    $ws.Range("Q4").Value = 2.381271987858333
    $ws.Range("R4").Value = 21.431447890725
    $ws.Range("S4").Value = 0.05787544222717098
    $ws.Range("T4").Value = 0.05787544222717096
    # Row 5
    $ws.Range("E5").Value = 3
    $ws.Range("F5").Value = 1
    $ws.Range("G5").Value = 1.419475
    $ws.Range("H5").Value = 4.258425
    $ws.Range("I5").Value = 0.1541931834006784
    $ws.Range("J5").Value = 0.1541931834006784
    $ws.Range("M5").Value = 0.7234496666666667
    $ws.Range("N5").Value = 2.170349
    $ws.Range("O5").Value = 0.1618662157416207
    $ws.Range("P5").Value = 0.1618662157416207
    $ws.Range("Q5").Value = 1.026918715591667
    $ws.Range("R5").Value = 9.242268440324999
    $ws.Range("S5").Value = 0.02495866709022151
    $ws.Range("T5").Value = 0.0249586670902215
    # Row 6
    $ws.Range("I6").Value = 0.3984988340349546
    $ws.Range("J6").Value = 0.3984988340349546
    $ws.Range("M6").Value = 1.716657
    $ws.Range("N6").Value = 5.149971
    $ws.Range("O6").Value = 0.3840886036988016
    $ws.Range("P6").Value = 0.3840886036988015
    $ws.Range("Q6").Value = 6.297572226631999
    $ws.Range("R6").Value = 56.678150039688
    $ws.Range("S6").Value = 0.1530588607400862
    $ws.Range("T6").Value = 0.1530588607400862
    # Row 7
    $ws.Range("I7").Value = 0.3984988340349546
    $ws.Range("J7").Value = 0.3984988340349546
    $ws.Range("O7").Value = 0.07870146593648156
    $ws.Range("P7").Value = 0.07870146593648154
    $ws.Range("S7").Value = 0.03136244241252959
    $ws.Range("T7").Value = 0.03136244241252959
    # Row 8
    $ws.Range("I8").Value = 0.3984988340349546
    $ws.Range("J8").Value = 0.3984988340349546
    $ws.Range("M8").Value = 1.677572333333333
    $ws.Range("N8").Value = 5.032717
    $ws.Range("O8").Value = 0.3753437146230962
    $ws.Range("P8").Value = 0.3753437146230962
    $ws.Range("Q8").Value = 6.15418976217511
    $ws.Range("R8").Value = 55.387707859576
    $ws.Range("S8").Value = 0.1495740326396526
    $ws.Range("T8").Value = 0.1495740326396525
    # Row 9
    $ws.Range("I9").Value = 0.3984988340349546
    $ws.Range("J9").Value = 0.3984988340349546
    $ws.Range("M9").Value = 0.7234496666666667
    $ws.Range("N9").Value = 2.170349
    $ws.Range("O9").Value = 0.1618662157416207
    $ws.Range("P9").Value = 0.1618662157416207
    $ws.Range("Q9").Value = 2.653981854363555
    $ws.Range("R9").Value = 23.885836689272
    $ws.Range("S9").Value = 0.06450349824268628
    $ws.Range("T9").Value = 0.06450349824268627
    # Row 10
    $ws.Range("G10").Value = 4.049549666666667
    $ws.Range("H10").Value = 12.148649
    $ws.Range("I10").Value = 0.4398900681184871
    $ws.Range("J10").Value = 0.439890068118487
    $ws.Range("M10").Value = 1.716657
    $ws.Range("N10").Value = 5.149971
    $ws.Range("O10").Value = 0.3840886036988016
    $ws.Range("P10").Value = 0.3840886036988015
    $ws.Range("Q10").Value = 6.951687782131
    $ws.Range("R10").Value = 62.56519003917899
    $ws.Range("S10").Value = 0.1689567620446004
    $ws.Range("T10").Value = 0.1689567620446004
    # Row 11
    $ws.Range("G11").Value = 4.049549666666667
    $ws.Range("H11").Value = 12.148649
    $ws.Range("I11").Value = 0.4398900681184871
    $ws.Range("J11").Value = 0.439890068118487
    $ws.Range("O11").Value = 0.07870146593648156
    $ws.Range("P11").Value = 0.07870146593648154
    $ws.Range("Q11").Value = 1.424431794949778
    $ws.Range("R11").Value = 12.819886154548
    $ws.Range("S11").Value = 0.03461999321182366
    $ws.Range("T11").Value = 0.03461999321182365
    # Row 12
    $ws.Range("G12").Value = 4.049549666666667
    $ws.Range("H12").Value = 12.148649
    $ws.Range("I12").Value = 0.4398900681184871
    $ws.Range("J12").Value = 0.439890068118487
    $ws.Range("M12").Value = 1.677572333333333
    $ws.Range("N12").Value = 5.032717
    $ws.Range("O12").Value = 0.3753437146230962
    $ws.Range("P12").Value = 0.3753437146230962
    $ws.Range("Q12").Value = 6.793412483259222
    $ws.Range("R12").Value = 61.14071234933299
    $ws.Range("S12").Value = 0.1651099721933998
    $ws.Range("T12").Value = 0.1651099721933997
    # Row 13
    $ws.Range("G13").Value = 4.049549666666667
    $ws.Range("H13").Value = 12.148649
    $ws.Range("I13").Value = 0.4398900681184871
    $ws.Range("J13").Value = 0.439890068118487
    $ws.Range("M13").Value = 0.7234496666666667
    $ws.Range("N13").Value = 2.170349
    $ws.Range("O13").Value = 0.1618662157416207
    $ws.Range("P13").Value = 0.1618662157416207
    $ws.Range("Q13").Value = 2.929645356500111
    $ws.Range("R13").Value = 26.36680820850099
    $ws.Range("S13").Value = 0.07120334066866328
    $ws.Range("T13").Value = 0.07120334066866325
    # Row 14
    $ws.Range("E14").Value = 1
    $ws.Range("F14").Value = 0.3333333333333333
    $ws.Range("G14").Value = 0.068288
    $ws.Range("H14").Value = 0.204864
    $ws.Range("I14").Value = 0.00741791444588001
    $ws.Range("J14").Value = 0.007417914445880009
    $ws.Range("M14").Value = 1.716657
    $ws.Range("N14").Value = 5.149971
    $ws.Range("O14").Value = 0.3840886036988016
    $ws.Range("P14").Value = 0.3840886036988015
    $ws.Range("Q14").Value = 0.117227073216
    $ws.Range("R14").Value = 1.055043658944
    $ws.Range("S14").Value = 0.002849136401875223
    $ws.Range("T14").Value = 0.002849136401875222
    # Row 15
    $ws.Range("E15").Value = 1
    $ws.Range("F15").Value = 0.3333333333333333
    $ws.Range("G15").Value = 0.068288
    $ws.Range("H15").Value = 0.204864
    $ws.Range("I15").Value = 0.00741791444588001
    $ws.Range("J15").Value = 0.007417914445880009
    $ws.Range("O15").Value = 0.07870146593648156
    $ws.Range("P15").Value = 0.07870146593648154
    $ws.Range("Q15").Value = 0.02402034952533334
    $ws.Range("R15").Value = 0.216183145728
    $ws.Range("S15").Value = 0.0005838007410821601
    $ws.Range("T15").Value = 0.0005838007410821599
    # Row 16
    $ws.Range("E16").Value = 1
    $ws.Range("F16").Value = 0.3333333333333333
    $ws.Range("G16").Value = 0.068288
    $ws.Range("H16").Value = 0.204864
    $ws.Range("I16").Value = 0.00741791444588001
    $ws.Range("J16").Value = 0.007417914445880009
    $ws.Range("M16").Value = 1.677572333333333
    $ws.Range("N16").Value = 5.032717
    $ws.Range("O16").Value = 0.3753437146230962
    $ws.Range("P16").Value = 0.3753437146230962
    $ws.Range("Q16").Value = 0.1145580594986667
    $ws.Range("R16").Value = 1.031022535488
    $ws.Range("S16").Value = 0.00278426756287293
    $ws.Range("T16").Value = 0.002784267562872929
    # Row 17
    $ws.Range("E17").Value = 1
    $ws.Range("F17").Value = 0.3333333333333333
    $ws.Range("G17").Value = 0.068288
    $ws.Range("H17").Value = 0.204864
    $ws.Range("I17").Value = 0.00741791444588001
    $ws.Range("J17").Value = 0.007417914445880009
    $ws.Range("M17").Value = 0.7234496666666667
    $ws.Range("N17").Value = 2.170349
    $ws.Range("O17").Value = 0.1618662157416207
    $ws.Range("P17").Value = 0.1618662157416207
    $ws.Range("Q17").Value = 0.04940293083733333
    $ws.Range("R17").Value = 0.4446263775359999
    $ws.Range("S17").Value = 0.001200709740049699
    $ws.Range("T17").Value = 0.001200709740049698
